# Journal de bord [Antoine] — add the 8 février 2016 entry
# (frmMotDePasse / frmPeriodes changes log row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entry row (row 5) ---------------------------------------
$ws.Range("A5").Value = "lundi 8 février 2016"

$ws.Range("B5").Value = "Mise à jour de toutes les formes, elles ne peuvent plus être redimensionnées. `nMise à jour de frmMotDePasse, confirmer que le mot de passe est conforme"

# Wrap the text and draw a medium border on the right edge of B5, like
# the rest of the table.
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Borders.Item(10).Weight = -4138

# Row height for the new row.
$ws.Rows.Item(5).RowHeight = 60

# --- Scroll the view down a bit and move the selection ----------------
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select() | Out-Null
